# corte de caja report
# Clear the date-range text from the merged A3:G3 header cell (was the
# shared string "Del viernes 24/junio/2022 al viernes 24/junio/2022"),
# leaving the cell blank but keeping its existing formatting/merge intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = ""

# Leave the cursor/selection where the author last clicked after editing.
[void]$ws.Range("C5").Select()
